$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: fill in End Time (D2), Total Time (E2), Job description (F2)
$ws.Range("D2").Value = 0.25
$ws.Range("E2").Value = "3 hours 20 minitues + seheri + other"
$ws.Range("F2").Value = "Black Box Testing Done"

# Row 3: fill in Date (B3), Start Time (C3), End Time (D3)
$ws.Range("B3").Value = "17-04-2022"
$ws.Range("C3").Value = 0.91666666666666663
$ws.Range("D3").Value = 0.9375

# Row 4: fill in Date (B4), Start Time (C4), End Time (D4)
$ws.Range("B4").Value = "17-04-2022"
$ws.Range("C4").Value = 0.029166666666666664
$ws.Range("D4").Value = 0.046527777777777779

# Update active sheet selection to E4
$ws.Activate()
$ws.Range("E4").Select()
